$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.710.12'
$ws.Range("E2").Value = '  -1.55%  '
$ws.Range("D3").Value = '3.391.94'
$ws.Range("E3").Value = '  -1.27%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = "'568.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.71%  '
$ws.Range("D6").Value = "'161.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.67%  '
$ws.Range("D8").Value = '3.394.57'
$ws.Range("E8").Value = '  -1.20%  '
$ws.Range("E9").Value = '  -5.54%  '
$ws.Range("E10").Value = '  +0.97%  '
$ws.Range("E11").Value = '  -2.02%  '
$ws.Range("D12").Value = "'0.422"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.62%  '
$ws.Range("D13").Value = '3.981.18'
$ws.Range("E13").Value = '  -1.19%  '
$ws.Range("E14").Value = '  +1.10%  '
$ws.Range("D15").Value = "'26.89"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.58%  '
$ws.Range("E16").Value = '  -2.96%  '
$ws.Range("D17").Value = '63.742.72'
$ws.Range("E17").Value = '  -1.57%  '
$ws.Range("D18").Value = '3.413.74'
$ws.Range("E18").Value = '  +0.12%  '
$ws.Range("D19").Value = "'6.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.03%  '
$ws.Range("D20").Value = "'13.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.17%  '
$ws.Range("D21").Value = "'375.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.21%  '
$ws.Range("D22").Value = "'7.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.60%  '
$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("D24").Value = "'70.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E25").Value = '  -6.03%  '
$ws.Range("E26").Value = '  -3.55%  '
$ws.Range("D27").Value = "'9.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.09%  '
$ws.Range("E28").Value = '  +0.30%  '
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("D30").Value = "'6.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.38%  '
$ws.Range("E31").Value = '  -6.01%  '
$ws.Range("E32").Value = '  -0.39%  '
$ws.Range("D33").Value = "'22.83"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.51%  '
$ws.Range("D34").Value = "'7.06"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.76%  '
$ws.Range("E35").Value = '  -4.62%  '
$ws.Range("D36").Value = "'159.41"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.28%  '
$ws.Range("D37").Value = "'0.856"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +10.08%  '
$ws.Range("E38").Value = '  -4.13%  '
$ws.Range("D39").Value = "'0.0724"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.59%  '
$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").Value = '2.763.55'
$ws.Range("E40").Value = '  -4.48%  '
$ws.Range("B41").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D41").Value = "'25.63"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.26%  '
$ws.Range("D42").Value = "'42.64"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.61%  '
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").Value = "'26.12"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.62%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = "'6.37"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.97%  '
$ws.Range("E45").Value = '  -3.46%  '
$ws.Range("D46").Value = "'0.0306"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.91%  '
$ws.Range("D47").Value = "'2.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +7.13%  '
$ws.Range("D48").Value = "'328.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.04%  '
$ws.Range("D49").Value = "'1.03"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.19%  '
$ws.Range("D50").Value = "'6.28"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.30%  '
$ws.Range("E51").Value = '  -2.96%  '
